# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E22) is re-ordered from descending
# (2210 .. 2204) to ascending (2204 .. 2210), and the "Valor Mora"
# values in F16 and F22 are swapped to match the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reverse the Periodo Mora period labels (E16:E22) so they read
# 2204, 2205, 2206, 2207, 2208, 2209, 2210 (ascending) instead of the
# previous descending order.
$ws.Range("E16").Value = "2204"
$ws.Range("E17").Value = "2205"
$ws.Range("E18").Value = "2206"
$ws.Range("E19").Value = "2207"
$ws.Range("E20").Value = "2208"
$ws.Range("E21").Value = "2209"
$ws.Range("E22").Value = "2210"

# Swap the Valor Mora amounts for the first and last rows so the
# amounts still line up with their (now reordered) period.
$ws.Range("F16").Value = 40000
$ws.Range("F22").Value = 36000
